# edit.ps1
#
# Applies the "prograse ppt and txt file" commit to
# PathFindingAndNevigationSystem.pptx:
#
#  1. Handout master footer date field text 06/01/2024 -> 14/04/2025.
#  2. On slide 2 (the presenter/author info slide):
#       - Remove the "Mr. Montaser Abdul Quader" title placeholder and the
#         "Lecturer, Department of CSE / Green University of Bangladesh"
#         subtitle placeholder entirely.
#       - Move + retext the "Student ID- ..." box to the new student id
#         232002280 at its new position.
#       - Move + retext the "Md. Mazharul Islam Shehab" box to read
#         "Rukonuzzaman Topu" at its new position.
#       - Remove the "Present to: " textbox and the straight connector
#         line next to it.
#       - Move the "Present by: " textbox to its new position.

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

function Remove-PlaceholderShapeCompletely($slide, $id) {
    # Deleting a placeholder shape (msoPlaceholder, Type=14) can make
    # PowerPoint immediately re-create an empty placeholder shell (from the
    # slide layout) under a brand new shape Id. Detect that by checking
    # whether a placeholder of the same PlaceholderFormat.Type still exists
    # afterwards, and delete it too (bounded loop just in case).
    $sh = Get-ShapeById $slide $id
    if ($sh -eq $null) { return }

    $phType = $null
    try { $phType = $sh.PlaceholderFormat.Type } catch {}

    $sh.Delete()

    for ($guard = 0; $guard -lt 4; $guard++) {
        $regenerated = $null
        for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
            $cand = $slide.Shapes.Item($i)
            if ($cand.Type -eq 14) {
                $candType = $null
                try { $candType = $cand.PlaceholderFormat.Type } catch {}
                if ($candType -eq $phType) {
                    $regenerated = $cand
                    break
                }
            }
        }
        if ($regenerated -eq $null) { break }
        $regenerated.Delete()
    }
}

function Set-ShapeEmuPosition($shape, $offXEmu, $offYEmu) {
    # Shape.Left / Shape.Top are expressed in points (1 pt = 12700 EMU) and
    # are stored internally with single-precision (float32) rounding, so a
    # naive `emu / 12700.0` can truncate down by 1 EMU once converted back.
    # Nudge the point value by a hair so it still rounds/truncates to the
    # exact target EMU after the float32 round-trip.
    $ptX = ($offXEmu / 12700.0) + 0.0000005
    $ptY = ($offYEmu / 12700.0) + 0.0000005
    $shape.Left = $ptX
    $shape.Top = $ptY
}

$p = $ppt.ActivePresentation

# --- Handout master: update the fixed footer date text --------------------
try {
    $hm = $p.HandoutMaster
    for ($i = 1; $i -le $hm.Shapes.Count; $i++) {
        $hsh = $hm.Shapes.Item($i)
        if ($hsh.HasTextFrame -and $hsh.TextFrame.HasText) {
            if ($hsh.TextFrame.TextRange.Text -eq "06/01/2024") {
                $hsh.TextFrame.TextRange.Text = "14/04/2025"
            }
        }
    }
} catch {}

# --- Slide 2: presenter / author info slide --------------------------------
$s = $p.Slides.Item(2)

# Remove the old "Mr. Montaser Abdul Quader" title and
# "Lecturer, Department of CSE / Green University of Bangladesh" subtitle.
Remove-PlaceholderShapeCompletely $s 548
Remove-PlaceholderShapeCompletely $s 549

# Update + reposition the student-id textbox.
$shStudent = Get-ShapeById $s 22
if ($shStudent -ne $null) {
    $full = $shStudent.TextFrame.TextRange.Text
    $old = "Student ID- 221002534"
    $idx = $full.IndexOf($old)
    if ($idx -ge 0) {
        $chars = $shStudent.TextFrame.TextRange.Characters($idx + 1, $old.Length)
        $chars.Text = "Student ID- 232002280"
    }
    Set-ShapeEmuPosition $shStudent 2948494 2282668
}

# Update + reposition the supervisor/presented-by-name textbox.
$shName = Get-ShapeById $s 27
if ($shName -ne $null) {
    $shName.TextFrame.TextRange.Text = "Rukonuzzaman Topu"
    Set-ShapeEmuPosition $shName 2639772 1688690
}

# Remove the vertical divider connector and the "Present to: " textbox.
$shConnector = Get-ShapeById $s 29
if ($shConnector -ne $null) { $shConnector.Delete() }

$shPresentTo = Get-ShapeById $s 30
if ($shPresentTo -ne $null) { $shPresentTo.Delete() }

# Reposition the "Present by: " textbox (text unchanged).
$shPresentBy = Get-ShapeById $s 31
if ($shPresentBy -ne $null) {
    Set-ShapeEmuPosition $shPresentBy 2855744 1385134
}
